$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# Row 2
Set-TextCell $ws.Range("D2") "329.16"
Set-TextCell $ws.Range("E2") "-0.95%"

# Row 3
Set-TextCell $ws.Range("D3") "41.28"
Set-TextCell $ws.Range("E3") "3.81%"

# Row 4
Set-TextCell $ws.Range("D4") "5.654"
Set-TextCell $ws.Range("E4") "-1.99%"

# Row 5
Set-TextCell $ws.Range("D5") "0.08314"
Set-TextCell $ws.Range("E5") "3.19%"

# Row 6
Set-TextCell $ws.Range("D6") "2.030"
Set-TextCell $ws.Range("E6") "2.54%"

# Row 7
Set-TextCell $ws.Range("D7") "8.770"
Set-TextCell $ws.Range("E7") "1.34%"

# Row 8
Set-TextCell $ws.Range("D8") "4.539"
Set-TextCell $ws.Range("E8") "0.92%"

# Row 9
Set-TextCell $ws.Range("D9") "2.946"
Set-TextCell $ws.Range("E9") "-1.76%"

# Row 10
Set-TextCell $ws.Range("D10") "0.9234"
Set-TextCell $ws.Range("E10") "0.06%"

# Row 11
Set-TextCell $ws.Range("D11") "0.1269"
Set-TextCell $ws.Range("E11") "-0.20%"

# Row 12
Set-TextCell $ws.Range("D12") "0.1959"
Set-TextCell $ws.Range("E12") "-0.17%"

# Row 13
Set-TextCell $ws.Range("D13") "0.09346"
Set-TextCell $ws.Range("E13") "0.66%"

# Row 14
Set-TextCell $ws.Range("D14") "0.03922"
Set-TextCell $ws.Range("E14") "9.98%"

# Row 15
Set-TextCell $ws.Range("D15") "0.1061"
Set-TextCell $ws.Range("E15") "0.91%"

# Row 16
Set-TextCell $ws.Range("D16") "0.001314"
Set-TextCell $ws.Range("E16") "0.52%"

# Row 17
Set-TextCell $ws.Range("D17") "0.006161"
Set-TextCell $ws.Range("E17") "-2.48%"

# Row 18
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextCell $ws.Range("D18") "3.438"
Set-TextCell $ws.Range("E18") "2.12%"

# Row 19
$ws.Range("B19").Value = "BitpandaEcosystemToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextCell $ws.Range("D19") "0.3534"
Set-TextCell $ws.Range("E19") "1.49%"

# Row 20
$ws.Range("B20").Value = "MCDex"
$ws.Range("C20").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws.Range("D20") "8.315"
Set-TextCell $ws.Range("E20") "-5.11%"

# Row 21
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextCell $ws.Range("D21") "0.1375"
Set-TextCell $ws.Range("E21") "1.11%"

# Row 22
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextCell $ws.Range("D22") "0.2446"
Set-TextCell $ws.Range("E22") "-11.33%"

# Row 23
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws.Range("D23") "0.04402"
Set-TextCell $ws.Range("E23") "-0.14%"

# Row 24
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextCell $ws.Range("D24") "0.001255"
Set-TextCell $ws.Range("E24") "-0.46%"

# Row 25
$ws.Range("B25").Value = "HotbitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextCell $ws.Range("D25") "0.004321"
Set-TextCell $ws.Range("E25") "-6.43%"

# Row 26
$ws.Range("B26").Value = "NitroEx"
$ws.Range("C26").Value = "https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx"
Set-TextCell $ws.Range("D26") "0.0001201"
Set-TextCell $ws.Range("E26") "0.95%"

# Row 27
$ws.Range("B27").Value = "UpBots"
$ws.Range("C27").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
Set-TextCell $ws.Range("D27") "0.0003049"
Set-TextCell $ws.Range("E27") "-23.64%"

# Row 39
Set-TextCell $ws.Range("D39") "0.02759"
Set-TextCell $ws.Range("E39") "11.35%"

# Row 40
Set-TextCell $ws.Range("D40") "0.05497"
Set-TextCell $ws.Range("E40") "-0.92%"

# Row 41
Set-TextCell $ws.Range("D41") "0.007784"
Set-TextCell $ws.Range("E41") "3.98%"

# Row 42
Set-TextCell $ws.Range("D42") "0.1421"
Set-TextCell $ws.Range("E42") "0.91%"

# Row 43
Set-TextCell $ws.Range("D43") "0.008945"
Set-TextCell $ws.Range("E43") "-9.99%"

# Row 44
Set-TextCell $ws.Range("E44") "1.65%"

# Row 45
Set-TextCell $ws.Range("D45") "0.01188"
Set-TextCell $ws.Range("E45") "11.92%"

# Row 46
Set-TextCell $ws.Range("D46") "0.00006979"
Set-TextCell $ws.Range("E46") "4.47%"

# Row 47
Set-TextCell $ws.Range("D47") "0.00000000751"
Set-TextCell $ws.Range("E47") "0.17%"

# Row 48
Set-TextCell $ws.Range("D48") "0.003184"
Set-TextCell $ws.Range("E48") "5.18%"

# Row 49
Set-TextCell $ws.Range("E49") "0.10%"

# Row 50
Set-TextCell $ws.Range("D50") "0.00002103"
Set-TextCell $ws.Range("E50") "0.17%"

# Row 51
Set-TextCell $ws.Range("D51") "0.0002003"
Set-TextCell $ws.Range("E51") "0.17%"
